# "parse opts from excel sheet, refine front-end app"
#
# The previously-blank "Sheet3" becomes a new "Config" sheet (moved to the
# front of the workbook) holding title/description/randomize* settings that
# the front-end now reads out of the workbook. The quiz data sheet ("Sheet1")
# stops being the tab that's selected/scrolled when the file is opened, since
# Config takes over as the active tab.

$wb = $excel.ActiveWorkbook

# Repurpose the existing empty "Sheet3" as the new "Config" sheet instead of
# inserting a brand new one, so it keeps its original sheetId and simply
# moves to the first tab position.
$cfg = $wb.Worksheets.Item("Sheet3")
$cfg.Name = "Config"
$cfg.Move($wb.Worksheets.Item("Sheet1"))
$cfg = $wb.Worksheets.Item("Config")

# Config key/value rows. B1 is written before A1 to match the order the
# values were originally entered in (and therefore the shared-string table
# order the workbook was saved with).
$cfg.Range("B1").Value = "My derpy test"
$cfg.Range("A1").Value = "title"

$cfg.Range("A2").Value = "description"
$cfg.Range("B2").Value = "This description comes from the excel doc"

$cfg.Range("A3").Value = "randomizeQuestions"
$cfg.Range("B3").Value = $true

$cfg.Range("A4").Value = "randomizeAnswers"
$cfg.Range("B4").Value = $true

# Size the two columns to fit their new contents.
$cfg.Columns.Item(1).EntireColumn.AutoFit()
$cfg.Columns.Item(2).EntireColumn.AutoFit()

# Config ends up the active tab/selection; the old quiz sheet ("Sheet1")
# naturally keeps its own last selection (C23) but is no longer the active
# tab or scrolled to row 16.
[void]$cfg.Range("C11").Select()
